$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("Unit Test Plan Preparation") - log 1 hr on day 9 (col L, previously blank)
# and correct day 12 (col O) from 2 hrs down to 1 hr
$ws.Range("L14").Value = 1
$ws.Range("O14").Value = 1

# Row 19 ("LLD Rework") - bump day 12 (col O) from 0.5 hr to 1.5 hrs
$ws.Range("O19").Value = 1.5

# Rows 24-31 - log hours on day 14 (col Q), previously blank
$ws.Range("Q24").Value = 1
$ws.Range("Q25").Value = 1
$ws.Range("Q26").Value = 1
$ws.Range("Q27").Value = 2
$ws.Range("Q28").Value = 2
$ws.Range("Q29").Value = 1
$ws.Range("Q30").Value = 1
$ws.Range("Q31").Value = 1

# Rows 32-33 - log hours on day 15 (col R), previously blank
$ws.Range("R32").Value = 3
$ws.Range("R33").Value = 1

# Leave the active selection where the editor's cursor ended up
$ws.Range("S18").Select()
